# Swap the data (columns B:AD) between each pair of rows, leaving the
# running index in column A untouched. This reflects the reordering of
# match records that occurred in the source data update.
#
# NOTE: we use .Value2 rather than .Value - in this host, reading/writing
# through the plain .Value property returns/accepts an opaque property
# descriptor instead of the actual cell contents, silently corrupting the
# data. .Value2 round-trips numbers/strings correctly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(124, 125),
    @(170, 171),
    @(194, 195)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AD$r1")
    $range2 = $ws.Range("B$r2`:AD$r2")

    $temp = $range1.Value2
    $range1.Value2 = $range2.Value2
    $range2.Value2 = $temp
}
